$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13700
$ws.Range("J2").Value = 20650
$ws.Range("L2").Formula = "=1040"
$ws.Range("B3").Value = 41660
$ws.Range("C3").Formula = "=650+750"
$ws.Range("F3").Value = 25020
$ws.Range("G3").Value = 600
$ws.Range("J3").Value = 6410
$ws.Range("B4").Value = 22250
$ws.Range("D4").Formula = "=600"
$ws.Range("F4").Value = 16370
$ws.Range("J4").Value = 21805
$ws.Range("B5").Value = 48350
$ws.Range("C5").Formula = "=1200+2670+2775+3620+3120"
$ws.Range("F5").Value = 20340
$ws.Range("G5").Formula = "=600"
$ws.Range("J5").Value = 8510
$ws.Range("B6").Value = 10470
$ws.Range("D6").Formula = "=1200"
$ws.Range("F6").Value = 11220
$ws.Range("H6").Formula = "=600"
$ws.Range("J6").Value = 4050
$ws.Range("K6").Formula = "=100+500"
$ws.Range("B7").Value = 10100
$ws.Range("C7").Formula = "=2000"
$ws.Range("F7").Value = 31860
$ws.Range("J7").Value = 13540
$ws.Range("B8").Value = 9615
$ws.Range("F8").Value = 9280
$ws.Range("J8").Value = 12145
$ws.Range("K8").Value = 760
$ws.Range("B9").Value = 23440
$ws.Range("D9").Value = 4370
$ws.Range("F9").Value = 14400
$ws.Range("J9").Value = 1220
$ws.Range("K9").Value = 650
$ws.Range("B10").Value = 20220
$ws.Range("F10").Value = 20400
$ws.Range("G10").Value = 835
$ws.Range("B11").Value = 20350
$ws.Range("F11").Value = 13120
$ws.Range("G11").Value = 7990
$ws.Range("B12").Value = 11400
$ws.Range("C12").Value = 1000
$ws.Range("F12").Value = 12660
$ws.Range("G12").Value = 4880
$ws.Range("J12").Value = 10830
$ws.Range("K12").Value = 1880
$ws.Range("B14").Value = 7150
$ws.Range("F14").Value = 9080
$ws.Range("J14").Value = 740
$ws.Range("B15").Value = 12630
$ws.Range("C15").Value = 1550
$ws.Range("F15").Value = 7100
$ws.Range("J15").Value = 45685
$ws.Range("L15").Value = 21400
$ws.Range("B16").Value = 12085
$ws.Range("F16").Value = 19850
$ws.Range("J16").Value = 14230
$ws.Range("K16").Value = 1105
$ws.Range("B17").Value = 34000
$ws.Range("D17").Value = 2800
$ws.Range("F17").Value = 14200
$ws.Range("G17").Value = 2040
$ws.Range("J17").Value = 13370
$ws.Range("K17").Value = 186
$ws.Range("L17").Value = 1640
$ws.Range("B18").Value = 27745
$ws.Range("C18").Value = 240
$ws.Range("F18").Value = 12500
$ws.Range("J18").Value = 11134
$ws.Range("L18").Value = 186
$ws.Range("B20").Value = 13050
$ws.Range("C20").Value = 100
$ws.Range("F20").Value = 12300
$ws.Range("J20").Value = 6695
$ws.Range("B21").Value = 28315
$ws.Range("F21").Value = 10850
$ws.Range("J21").Value = 9250
$ws.Range("K21").Value = 450
$ws.Range("B22").Value = 11335
$ws.Range("F22").Value = 7460
$ws.Range("J22").Value = 12750
$ws.Range("K22").Value = 1920
$ws.Range("B23").Value = 15505
$ws.Range("F23").Value = 19800
$ws.Range("J23").Value = 8630
$ws.Range("B24").Value = 5935
$ws.Range("C24").Value = 320
$ws.Range("J24").Value = 3700
$ws.Range("B25").Value = 15590
$ws.Range("J25").Value = 8460
$ws.Range("B26").Value = 11940
$ws.Range("F26").Value = 26290
$ws.Range("J26").Value = 15440
$ws.Range("B27").Value = 19720
$ws.Range("C27").Value = 5820
$ws.Range("F27").Value = 13170
$ws.Range("B28").Value = 26660
$ws.Range("C28").Value = 3160
$ws.Range("F28").Value = 18890
$ws.Range("J28").Value = 13800
$ws.Range("L28").Value = 320
$ws.Range("B29").Value = 9850
$ws.Range("D29").Value = 100
$ws.Range("F29").Value = 15920
$ws.Range("J29").Value = 12840
$ws.Range("F30").Value = 19000
$ws.Range("J30").Value = 46300
$ws.Range("K30").Value = 13900
$ws.Range("B31").Formula = "=21800+4220"
$ws.Range("C31").Value = 5200
$ws.Range("D31").Value = 1755
$ws.Range("F31").Value = 14210

# Style fix-ups: a couple of cells change their cached style index
# (F3 goes from style "6" to style "5"; K17 is a brand-new cell that
# picks up style "6"). Borrow the format from an untouched donor cell
# that already carries the desired style so the xf index matches.
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H3").Copy()
$ws.Range("K17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Column width adjustments (best-effort; engine quantizes to 1/6 px grid)
$ws.Columns.Item(5).ColumnWidth = 13.736979166666666
$ws.Columns.Item(9).ColumnWidth = 16.307291666666668
$ws.Columns.Item(10).ColumnWidth = 10.022135416666666
$ws.Columns.Item(11).ColumnWidth = 10.022135416666666
$ws.Columns.Item(12).ColumnWidth = 12.877604166666666

# Final selected cell (matches last active cell recorded in the sheet view)
[void]$ws.Range("F12").Select()
